$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used with PasteSpecial to copy only cell formatting.
$xlPasteFormats = -4122

# --- Row 4: new "11 or 12" key entry (green fill, matching the "Cog" block style) ---
$ws.Range("A10").Copy() | Out-Null
$ws.Range("E4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E4").Value = "11 or 12"

# --- Row 5: new "33 or 34" key entry (same green fill) ---
$ws.Range("A10").Copy() | Out-Null
$ws.Range("E5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E5").Value = "33 or 34"

# --- Row 7: F7 switches from the green "11 & 12" entry to a gold "8 & 19" entry ---
$ws.Range("A6").Copy() | Out-Null
$ws.Range("F7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F7").Value = "8 & 19"

# --- Row 8: new "2 & 14" key entry (blue-tint fill, matching F4:F6) ---
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F8").Value = "2 & 14"

# --- Row 9: new "7 & 19" key entry (gold fill, matching A6:C9 / F7) ---
$ws.Range("A6").Copy() | Out-Null
$ws.Range("F9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F9").Value = "7 & 19"

# --- Row 6: new "9 or 10" key entry (same green fill) ---
$ws.Range("A10").Copy() | Out-Null
$ws.Range("E6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E6").Value = "9 or 10"

$excel.CutCopyMode = $false

# Match the saved selection / active cell from the edited workbook.
$ws.Range("E6").Select()
